# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" footer timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 06:03"

# Swap the ordering of Curazao/Dominica (rows 198-199) and
# Seychelles/Montserrat (rows 205-206), along with their stats.
$ws.Range("A198").Value = "Curazao"
$ws.Range("A199").Value = "Dominica"

$ws.Range("A205").Value = "Seychelles"
$ws.Range("A206").Value = "Montserrat"

# Update statistics (Casos activos = D, Casos criticos = F, Muertes hoy = G, Muertes = H)

# Row 18 - India
$ws.Range("D18").Value = 14183
$ws.Range("E18").Value = 33523
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 1694

# Row 198 - Curazao
$ws.Range("D198").Value = 13
$ws.Range("H198").Value = 1

# Row 199 - Dominica
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 0

# Row 205 - Seychelles
$ws.Range("D205").Value = 8
$ws.Range("F205").Value = 0
$ws.Range("H205").Value = 0

# Row 206 - Montserrat
$ws.Range("D206").Value = 7
$ws.Range("F206").Value = 1
$ws.Range("H206").Value = 1
